$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (D) and 1h volume-change (E) figures to the latest
# scrape snapshot. Cells store these as literal text (e.g. "309.32", "-0.01%"),
# so NumberFormat is forced to Text ("@") before assignment to stop Excel from
# reinterpreting the numeric- and percent-looking strings as real numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "309.32"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.01%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "41.09"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "0.15%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.193"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.51%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07681"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.65%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.302"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.26%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.684"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "4.23%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9153"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "1.10%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.425"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-2.37%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1240"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "11.49%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1827"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "1.80%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09148"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "1.11%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.04186"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-1.64%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.1051"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.16%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001309"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "4.97%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005778"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "1.88%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.344"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.10%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.06%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.462"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "12.08%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1374"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.75%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.2824"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "4.34%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.04009"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-0.50%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.001269"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.00%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.004097"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.10%"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.45%"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02531"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "5.02%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05311"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "1.29%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007854"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "1.25%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1309"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "0.65%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.006647"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.001877"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-3.67%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007411"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-12.15%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3064"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-8.24%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006791"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-1.48%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.45%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.2334"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "293.08%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002108"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.45%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002008"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.45%"
